# Add a "Save" column (H) to the s_vals sheet, mirroring the existing
# header style used by the other stat columns (e.g. G1 "sum") and filling
# in the per-row Save flag (0/1) for rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the neighboring header cell (G1) onto H1, then
# overwrite its text with the new header label "Save".
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Per-row "Save" values for rows 2-10.
$saveValues = @(0, 0, 0, 1, 0, 0, 0, 0, 1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
